$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current (pre-edit) values for columns D, L, M, N, O, P, Q, R, S, T across rows 2-20
# so the row permutation below reads consistent "before" data regardless of write order.
$snapshot = @{}
for ($r = 2; $r -le 20; $r++) {
    $row = @{}
    $row["D"] = $ws.Cells.Item($r, 4).Value2
    $row["L"] = $ws.Cells.Item($r, 12).Value()
    $row["M"] = $ws.Cells.Item($r, 13).Value2
    $row["N"] = $ws.Cells.Item($r, 14).Value2
    $row["O"] = $ws.Cells.Item($r, 15).Value2
    $row["P"] = $ws.Cells.Item($r, 16).Value2
    $row["Q"] = $ws.Cells.Item($r, 17).Value()
    $row["R"] = $ws.Cells.Item($r, 18).Value()
    $row["S"] = $ws.Cells.Item($r, 19).Value2
    $row["T"] = $ws.Cells.Item($r, 20).Value2
    $snapshot[$r] = $row
}

# Destination row <- source row (original row number), per the weekly re-shuffle in the commit
$rowSource = @{
    2 = 3
    3 = 4
    4 = 18
    5 = 10
    6 = 13
    7 = 14
    8 = 11
    9 = 8
    10 = 17
    11 = 6
    12 = 2
    13 = 16
    14 = 20
    15 = 5
    16 = 19
    17 = 7
    18 = 9
    19 = 15
    20 = 12
}

foreach ($destRow in 2..20) {
    $srcRow = $rowSource[$destRow]
    $data = $snapshot[$srcRow]
    $ws.Cells.Item($destRow, 4).Value2 = $data["D"]
    $ws.Cells.Item($destRow, 12).Value = $data["L"]
    $ws.Cells.Item($destRow, 13).Value2 = $data["M"]
    $ws.Cells.Item($destRow, 14).Value2 = $data["N"]
    $ws.Cells.Item($destRow, 15).Value2 = $data["O"]
    $ws.Cells.Item($destRow, 16).Value2 = $data["P"]
    $ws.Cells.Item($destRow, 17).Value = $data["Q"]
    $ws.Cells.Item($destRow, 18).Value = $data["R"]
    $ws.Cells.Item($destRow, 19).Value2 = $data["S"]
    $ws.Cells.Item($destRow, 20).Value2 = $data["T"]
}

# Row 20's "Precio $/Kg" (S) keeps its original reported value (2500) rather than
# following the rest of the row's source data - matches the published weekly figures.
$ws.Cells.Item(20, 19).Value2 = $snapshot[20]["S"]

Write-Host "Done reshuffling rows 2-20"
